# Refresh market-price-derived Leve profit columns (H,I,J,K,L,M,N) across sheets.
# Mirrors a scheduled market-data re-pull: most rows only change cached values,
# a few rows gain/lose the NQ or HQ profit cell depending on sign flips.
$wb = $excel.ActiveWorkbook

# ALC row 61: Not Taking No for an Answer / Mega-Potion of Strength
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H61").Value = 5735.3335
$ws.Range("J61").Value = 17017
$ws.Range("L61").Value = 51051
$ws.Range("N61").Value = -51395

# ALC row 114: Conserving Combat / Bluespirit Codex
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H114").Value = 39993
$ws.Range("J114").Value = 39993
$ws.Range("L114").Value = 39993
$ws.Range("N114").Value = -48671

# ARM row 37: Get Shirty / Steel Chainmail
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 14685.667
$ws.Range("I37").Value = 5000
$ws.Range("J37").Value = 19528.5
$ws.Range("K37").Value = 5000
$ws.Range("L37").Value = 19528.5
$ws.Range("M37").Value = -4727
$ws.Range("N37").Value = -20074.5

# ARM row 49: I've Got You under My Skin / Steel-plated Caligae
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()

# ARM row 103: Sweeping the Legs / Doman Steel Greaves of Striking
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H103").Value = 30958
$ws.Range("J103").Value = 30958
$ws.Range("L103").Value = 30958
$ws.Range("N103").Value = -33302

# ARM row 122: Haste for High Durium / High Durium Nugget
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2506.2646
$ws.Range("I122").Value = 2502.3794
$ws.Range("K122").Value = 7507.138199999999
$ws.Range("M122").Value = -5057.138199999999

# ARM row 138: Don't Ask about the Rivets / Titanium Gold Helm of Casting
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H138").Value = 31000
$ws.Range("J138").Value = 31000
$ws.Range("L138").Value = 31000
$ws.Range("N138").Value = -41280

# BSM row 32: Time to Upgrade / Iron Ornamental Hammer
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H32").Value = 7500
$ws.Range("I32").Value = 5000
$ws.Range("J32").Value = 10000
$ws.Range("K32").Value = 5000
$ws.Range("L32").Value = 10000
$ws.Range("M32").Value = -4616
$ws.Range("N32").Value = -10768

# BSM row 102: Renting Mortality / Doman Steel Mortar
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H102").Value = 2778
$ws.Range("I102").Value = 2778
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2778
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("M102").Value = 467

# BSM row 134: Ruthenium Supremium / Ruthenium Ingot
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 17882.508
$ws.Range("I134").Value = 21071.424
$ws.Range("J134").Value = 6827.6
$ws.Range("K134").Value = 63214.272
$ws.Range("L134").Value = 20482.8
$ws.Range("M134").Value = -60679.272
$ws.Range("N134").Value = -25552.8

# CRP row 7: Gridania's Got Talent / Maple Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 85.86667
$ws.Range("I7").Value = 67.333336
$ws.Range("J7").Value = 160
$ws.Range("K7").Value = 67.333336
$ws.Range("L7").Value = 160
$ws.Range("M7").Value = 45.666664
$ws.Range("N7").Value = -386

# CRP row 9: Shields for the Serpents / Round Shield
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 18487.334
$ws.Range("J9").Value = 18487.334
$ws.Range("L9").Value = 18487.334
$ws.Range("N9").Value = -18823.334

# CRP row 22: Driving Up the Wall / Elm Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 669.2727
$ws.Range("I22").Value = 345.18182
$ws.Range("J22").Value = 993.36365
$ws.Range("K22").Value = 345.18182
$ws.Range("L22").Value = 993.36365
$ws.Range("M22").Value = 4.818179999999984
$ws.Range("N22").Value = -1693.36365

# CRP row 58: You Do the Heavy Lifting / Mahogany Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2212.4666
$ws.Range("I58").Value = 1720.3529
$ws.Range("J58").Value = 2856
$ws.Range("K58").Value = 1720.3529
$ws.Range("L58").Value = 2856
$ws.Range("M58").Value = -1517.3529
$ws.Range("N58").Value = -3262

# CRP row 107: Built to Last / White Oak Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 306.94446
$ws.Range("I107").Value = 241.66667
$ws.Range("K107").Value = 241.66667
$ws.Range("M107").Value = 1678.33333

# CRP row 132: Hull Lotta Damage / Ginseng Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2172.7908
$ws.Range("I132").Value = 998.1579
$ws.Range("J132").Value = 3102.7083
$ws.Range("K132").Value = 2994.4737
$ws.Range("L132").Value = 9308.124899999999
$ws.Range("M132").Value = -464.4737
$ws.Range("N132").Value = -14368.1249

# CRP row 136: Turali Quality / Dark Mahogany Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2212.4666
$ws.Range("I136").Value = 1720.3529
$ws.Range("J136").Value = 2856
$ws.Range("K136").Value = 5161.0587
$ws.Range("L136").Value = 8568
$ws.Range("M136").Value = -2611.0587
$ws.Range("N136").Value = -13668

# CUL row 48: Rise and Dine / Cheese Souffle
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H48").Value = 3316.3
$ws.Range("I48").Value = 799.6667
$ws.Range("J48").Value = 4394.857
$ws.Range("K48").Value = 2399.0001
$ws.Range("L48").Value = 13184.571
$ws.Range("M48").Value = -2149.0001
$ws.Range("N48").Value = -13684.571

# CUL row 81: It Goes Down Smoothly / Frozen Spirits
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 5083.3335
$ws.Range("I81").Value = 1200
$ws.Range("J81").Value = 5860
$ws.Range("K81").Value = 3600
$ws.Range("L81").Value = 17580
$ws.Range("M81").Value = -2477
$ws.Range("N81").Value = -19826

# CUL row 84: Quenching the Flame (L) / Frozen Spirits
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H84").Value = 5083.3335
$ws.Range("I84").Value = 1200
$ws.Range("J84").Value = 5860
$ws.Range("K84").Value = 10800
$ws.Range("L84").Value = 52740
$ws.Range("M84").Value = -5184
$ws.Range("N84").Value = -63972

# CUL row 133: Friends Are Food / Boiled Alpaca Steak
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 2719.6191
$ws.Range("I133").Value = 2058.8462
$ws.Range("J133").Value = 3793.375
$ws.Range("K133").Value = 6176.5386
$ws.Range("L133").Value = 11380.125
$ws.Range("M133").Value = -1116.5386
$ws.Range("N133").Value = -21500.125

# LTW row 4: Sole Traders / Leather Duckbills
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 3276.6667
$ws.Range("I4").Value = 2615
$ws.Range("J4").Value = 4600
$ws.Range("K4").Value = 2615
$ws.Range("L4").Value = 4600
$ws.Range("M4").Value = -2502
$ws.Range("N4").Value = -4826

# LTW row 7: Tan Before the Ban / Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2144.24
$ws.Range("I7").Value = 2035.7858
$ws.Range("K7").Value = 2035.7858
$ws.Range("M7").Value = -1923.7858

# LTW row 22: Skin off Their Backs / Aldgoat Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 403.6
$ws.Range("I22").Value = 317.5
$ws.Range("J22").Value = 748
$ws.Range("K22").Value = 317.5
$ws.Range("L22").Value = 748
$ws.Range("M22").Value = -22.5
$ws.Range("N22").Value = -1338

# LTW row 27: Fire and Hide / Aldgoat Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 403.6
$ws.Range("I27").Value = 317.5
$ws.Range("J27").Value = 748
$ws.Range("K27").Value = 317.5
$ws.Range("L27").Value = 748
$ws.Range("M27").Value = -210.5
$ws.Range("N27").Value = -962

# LTW row 28: My Sole to Take / Padded Leather Duckbills
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H28").Value = 3276.6667
$ws.Range("I28").Value = 2615
$ws.Range("J28").Value = 4600
$ws.Range("K28").Value = 2615
$ws.Range("L28").Value = 4600
$ws.Range("M28").Value = -2383
$ws.Range("N28").Value = -5064

# LTW row 37: Quicker than Sand / Padded Leather Duckbills
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H37").Value = 3276.6667
$ws.Range("I37").Value = 2615
$ws.Range("J37").Value = 4600
$ws.Range("K37").Value = 2615
$ws.Range("L37").Value = 4600
$ws.Range("M37").Value = -2508
$ws.Range("N37").Value = -4814

# LTW row 46: Supply Side Logic / Boar Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 762.9032
$ws.Range("I46").Value = 720
$ws.Range("J46").Value = 798.2353000000001
$ws.Range("K46").Value = 720
$ws.Range("L46").Value = 798.2353000000001
$ws.Range("M46").Value = -532
$ws.Range("N46").Value = -1174.2353

# LTW row 55: It's Not a Job, It's a Calling / Peiste Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 316.25
$ws.Range("J55").Value = 276.3
$ws.Range("L55").Value = 276.3
$ws.Range("N55").Value = -622.3

# LTW row 82: Trainin' the Neck / Dragon Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2092.3784
$ws.Range("I82").Value = 1982.4482
$ws.Range("J82").Value = 2490.875
$ws.Range("K82").Value = 1982.4482
$ws.Range("L82").Value = 2490.875
$ws.Range("M82").Value = -1621.4482
$ws.Range("N82").Value = -3212.875

# LTW row 85: Training Is Only Skintight (L) / Dragon Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 2092.3784
$ws.Range("I85").Value = 1982.4482
$ws.Range("J85").Value = 2490.875
$ws.Range("K85").Value = 1982.4482
$ws.Range("L85").Value = 2490.875
$ws.Range("M85").Value = -734.4482
$ws.Range("N85").Value = -4986.875

# LTW row 122: Hell on Leather / Gaja Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2924.4119
$ws.Range("I122").Value = 2781
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 8343
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -5893
$ws.Range("N122").Value = -16900

# LTW row 126: Battered Books / Saiga Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2144.24
$ws.Range("I126").Value = 2035.7858
$ws.Range("K126").Value = 6107.357400000001
$ws.Range("M126").Value = -3637.357400000001

# LTW row 132: Tenets of Tanning / Silver Lobo Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 10635.667
$ws.Range("I132").Value = 3967.3333
$ws.Range("J132").Value = 12858.444
$ws.Range("K132").Value = 11901.9999
$ws.Range("L132").Value = 38575.33199999999
$ws.Range("M132").Value = -9371.999899999999
$ws.Range("N132").Value = -43635.33199999999

# WVR row 113: A Tender Table / Pixie Floss
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 354.72726
$ws.Range("I113").Value = 330.05
$ws.Range("J113").Value = 601.5
$ws.Range("K113").Value = 990.1500000000001
$ws.Range("L113").Value = 1804.5
$ws.Range("M113").Value = 1179.85
$ws.Range("N113").Value = -6144.5
